$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Version: 1.0 -> 1.2.5
$ws.Range("D2").Value = "1.2.5"

# Precondition text correction (appears for each of the 5 test cases)
$preconditionNew = "O usuário devidamente autenticado e na tela inicial do sistema."
foreach ($cell in @("B8", "B16", "B23", "B30", "B37")) {
    $ws.Range($cell).Value = $preconditionNew
}

# TC1 actual-result wording fix
$ws.Range("B10").Value = "Beneficiário Acessa o caso de uso através do menu."

# Add trailing periods to expected results text
$ws.Range("D18").Value = "SYSTEM Apresenta a tela de Cancelar Solicitação de Diária."
$ws.Range("D25").Value = "SYSTEM Apresenta a tela de Analisar Prestação de Contas."
$ws.Range("D32").Value = "SYSTEM Apresenta a tela de Detalhar Diárias."
